$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.149.36"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "2.419.48"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'558.58"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").Value = "'143.42"
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("D9").Value = "2.416.39"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("E13").Value = "  +0.94%  "
$ws.Range("D14").Value = "'26.26"
$ws.Range("E14").Value = "  +4.27%  "
$ws.Range("E15").Value = "  +5.40%  "
$ws.Range("D16").Value = "2.848.50"
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("D17").Value = "61.994.30"
$ws.Range("E17").Value = "  +1.36%  "
$ws.Range("D18").Value = "2.415.46"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "'11.19"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").Value = "'4.20"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").Value = "'324.47"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "'6.77"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'65.48"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").Value = "'9.01"
$ws.Range("E26").Value = "  +6.40%  "
$ws.Range("D27").Value = "'596.42"
$ws.Range("E27").Value = "  +17.29%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").Value = "2.521.62"
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("D30").Value = "0.0₃0941"
$ws.Range("E30").Value = "  +5.51%  "
$ws.Range("D31").Value = "'8.32"
$ws.Range("E31").Value = "  +1.56%  "
$ws.Range("E32").Value = "  +5.09%  "
$ws.Range("D33").Value = "'0.149"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("E35").Value = "  +2.37%  "
$ws.Range("D36").Value = "'5.73"
$ws.Range("E36").Value = "  +5.65%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'4.80"
$ws.Range("E38").Value = "  +2.17%  "
$ws.Range("E39").Value = "  +1.21%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").Value = "'151.45"
$ws.Range("E40").Value = "  +3.27%  "
$ws.Range("B41").Value = "EthereumClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D41").Value = "'18.72"
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("D42").Value = "'1.82"
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "'2.37"
$ws.Range("E44").Value = "  +12.38%  "
$ws.Range("D45").Value = "'151.14"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").Value = "'0.0541"
$ws.Range("E47").Value = "  +3.41%  "
$ws.Range("E48").Value = "  +4.49%  "
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").Value = "'0.0230"
$ws.Range("E51").Value = "  +2.11%  "
